$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.573.24"
$ws.Range("E2").Value = "  +2.10%  "
$ws.Range("D3").Value = "1.676.81"
$ws.Range("E3").Value = "  +2.67%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'220.21"
$ws.Range("E5").Value = "  +2.66%  "
$ws.Range("D6").Value = "'0.531"
$ws.Range("E6").Value = "  +2.57%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'30.04"
$ws.Range("E8").Value = "  +4.76%  "
$ws.Range("D9").Value = "'0.265"
$ws.Range("E9").Value = "  +2.49%  "
$ws.Range("D10").Value = "'0.0636"
$ws.Range("E10").Value = "  +4.47%  "
$ws.Range("E11").Value = "  -0.66%  "
$ws.Range("D12").Value = "1.919.94"
$ws.Range("E12").Value = "  +2.83%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "'0.617"
$ws.Range("E13").Value = "  +9.43%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'10.27"
$ws.Range("E14").Value = "  +10.26%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.673.75"
$ws.Range("E15").Value = "  +2.51%  "
$ws.Range("E16").Value = "  +3.31%  "
$ws.Range("D17").Value = "30.592.06"
$ws.Range("E17").Value = "  +2.18%  "
$ws.Range("D18").Value = "'66.46"
$ws.Range("E18").Value = "  +3.66%  "
$ws.Range("D19").Value = "'245.44"
$ws.Range("E19").Value = "  +0.96%  "
$ws.Range("D20").Value = "0.0₃0724"
$ws.Range("E20").Value = "  +3.23%  "
$ws.Range("D21").Value = "'1.00"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "'4.27"
$ws.Range("E22").Value = "  +3.58%  "
$ws.Range("D23").Value = "'10.06"
$ws.Range("E23").Value = "  +2.15%  "
$ws.Range("E24").Value = "  +1.47%  "
$ws.Range("D25").Value = "'158.03"
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("D26").Value = "'15.97"
$ws.Range("E26").Value = "  +2.82%  "
$ws.Range("E27").Value = "  +2.60%  "
$ws.Range("D28").Value = "'6.70"
$ws.Range("E28").Value = "  +1.58%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("E30").Value = "  +2.28%  "
$ws.Range("E31").Value = "  +3.17%  "
$ws.Range("E32").Value = "  +3.33%  "
$ws.Range("D33").Value = "1.512.66"
$ws.Range("E33").Value = "  +6.21%  "
$ws.Range("E34").Value = "  +4.11%  "
$ws.Range("D35").Value = "'1.77"
$ws.Range("E35").Value = "  +7.75%  "
$ws.Range("B36").Value = "Aave"
$ws.Range("C36").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D36").Value = "'84.36"
$ws.Range("E36").Value = "  +11.30%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "'1.03"
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("D38").Value = "'0.602"
$ws.Range("E38").Value = "  +9.02%  "
$ws.Range("E39").Value = "  +5.56%  "
$ws.Range("E41").Value = "  +0.48%  "
$ws.Range("D42").Value = "'0.841"
$ws.Range("E42").Value = "  +1.50%  "
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("D44").Value = "'0.0499"
$ws.Range("E44").Value = "  +1.89%  "
$ws.Range("E45").Value = "  +1.02%  "
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D48").Value = "'51.55"
$ws.Range("E48").Value = "  -3.01%  "
$ws.Range("D49").Value = "1.813.79"
$ws.Range("E49").Value = "  +2.13%  "
$ws.Range("D50").Value = "'95.02"
$ws.Range("E50").Value = "  +6.40%  "
$ws.Range("D51").Value = "0.0₆0113"
$ws.Range("E51").Value = "  +0.83%  "
